$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.849.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.64%  "

$ws.Range("D3").Value = "'1.636.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.62%  "

$ws.Range("D5").Value = "'214.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.58%  "

$ws.Range("D6").Value = "'0.5018"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("D7").Value = "'1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.67%  "

$ws.Range("D8").Value = "'0.2558"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.91%  "

$ws.Range("D9").Value = "'0.06357"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.07%  "

$ws.Range("D10").Value = "'19.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.87%  "

$ws.Range("D11").Value = "'0.07782"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("D12").Value = "'1.651.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.90%  "

$ws.Range("D13").Value = "'4.237"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.68%  "

$ws.Range("D14").Value = "'1.864.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("D15").Value = "'0.5393"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.31%  "

$ws.Range("D16").Value = "'0.0₅7848"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.62%  "

$ws.Range("D17").Value = "'64.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.07%  "

$ws.Range("D18").Value = "'25.943.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").Value = "'1.005"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.91%  "

$ws.Range("D20").Value = "'194.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.23%  "

$ws.Range("D21").Value = "'4.350"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.67%  "

$ws.Range("D22").Value = "'9.856"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.67%  "

$ws.Range("D23").Value = "'5.940"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.67%  "

$ws.Range("D24").Value = "'1.010"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.20%  "

$ws.Range("D25").Value = "'1.891"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.71%  "

$ws.Range("D26").Value = "'139.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.29%  "

$ws.Range("D27").Value = "'0.1126"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.39%  "

$ws.Range("D28").Value = "'6.775"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.68%  "

$ws.Range("D29").Value = "'15.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.50%  "

$ws.Range("D30").Value = "'1.237"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.17%  "

$ws.Range("D31").Value = "'0.04837"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.45%  "

$ws.Range("D32").Value = "'3.227"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.43%  "

$ws.Range("D33").Value = "'3.155"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.48%  "

$ws.Range("D34").Value = "'1.523"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.31%  "

$ws.Range("D35").Value = "'2.364"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.07%  "

$ws.Range("D36").Value = "'2.598"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.57%  "

$ws.Range("D37").Value = "'0.8811"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.23%  "

$ws.Range("D38").Value = "'1.124.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.38%  "

$ws.Range("D39").Value = "'0.5482"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.32%  "

$ws.Range("D40").Value = "'0.01563"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.11%  "

$ws.Range("D41").Value = "'1.004"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.77%  "

$ws.Range("D42").Value = "'5.640"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.37%  "

$ws.Range("D43").Value = "'0.8110"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.98%  "

$ws.Range("D44").Value = "'99.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.67%  "

$ws.Range("D45").Value = "'1.775.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.24%  "

$ws.Range("E46").Value = "  +4.88%  "

$ws.Range("D47").Value = "'0.4517"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.53%  "

$ws.Range("D48").Value = "'1.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.93%  "

$ws.Range("D49").Value = "'55.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.12%  "

$ws.Range("D50").Value = "'0.05036"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.09%  "

$ws.Range("D51").Value = "'1.008"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.14%  "
